# Convention change to support multi-axle vehicles:
# rename the front/rear axle parameter labels ("sAxleF"/"sAxleR") to the
# generic, ordinal labels "sAxle1"/"sAxle2" on every Body data sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sedan_HambaLG", "Sedan_Hamba", "Bus_Makhulu")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A5").Value = "sAxle1"
    $ws.Range("A6").Value = "sAxle2"
}

# Restore/update the per-sheet window selection state (mirrors what Excel
# records after a user reviewed Sedan_Hamba and Bus_Makhulu before coming
# back to rest on Sedan_HambaLG, cell A2).
$ws2 = $wb.Worksheets.Item("Sedan_Hamba")
$ws2.Activate()
$ws2.Range("A17").Select()

$ws3 = $wb.Worksheets.Item("Bus_Makhulu")
$ws3.Activate()
$ws3.Range("A2").Select()

$ws1 = $wb.Worksheets.Item("Sedan_HambaLG")
$ws1.Activate()
$ws1.Range("A2").Select()
